$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.45204894900344
$ws.Range("C2").Value = 9.080804078643411
$ws.Range("D2").Value = 6.018171943034317
$ws.Range("E2").Value = 10.68549349764371
$ws.Range("G2").Value = 50.79808072400368
$ws.Range("H2").Value = 19.31155231018177
$ws.Range("L2").Value = 9.968694094866992
$ws.Range("M2").Value = 17.19136353697392
$ws.Range("N2").Value = 20.15751439965944
$ws.Range("B3").Value = 20.03759347642024
$ws.Range("C3").Value = 8.606517243429112
$ws.Range("D3").Value = 5.906311185156739
$ws.Range("E3").Value = 10.68769380527232
$ws.Range("G3").Value = 50.37921649700017
$ws.Range("H3").Value = 19.30643426136705
$ws.Range("L3").Value = 9.980167639285536
$ws.Range("M3").Value = 17.11871649359626
$ws.Range("N3").Value = 20.22655325915195
$ws.Range("B4").Value = 19.78479836241165
$ws.Range("C4").Value = 8.299909035692977
$ws.Range("D4").Value = 5.838510401150267
$ws.Range("E4").Value = 10.68923898951661
$ws.Range("G4").Value = 50.13726530273155
$ws.Range("H4").Value = 19.30779418883041
$ws.Range("L4").Value = 9.988666276070864
$ws.Range("M4").Value = 17.07782932746373
$ws.Range("N4").Value = 20.27092681785572
$ws.Range("B5").Value = 19.6823605688284
$ws.Range("C5").Value = 8.171151691034046
$ws.Range("D5").Value = 5.811143086914735
$ws.Range("E5").Value = 10.68991747830264
$ws.Range("G5").Value = 50.04258518485096
$ws.Range("H5").Value = 19.30947872738004
$ws.Range("L5").Value = 9.99249515962582
$ws.Range("M5").Value = 17.06211428875167
$ws.Range("N5").Value = 20.28950918629425
$ws.Range("B6").Value = 19.66539034557965
$ws.Range("C6").Value = 8.149543370625448
$ws.Range("D6").Value = 5.806615765692714
$ws.Range("E6").Value = 10.69003308865582
$ws.Range("G6").Value = 50.02710255321645
$ws.Range("H6").Value = 19.30982664819498
$ws.Range("L6").Value = 9.993153025500126
$ws.Range("M6").Value = 17.05956232434637
$ws.Range("N6").Value = 20.29262499072463
$ws.Range("B7").Value = 19.78341430284123
$ws.Range("C7").Value = 8.298187916509194
$ws.Range("D7").Value = 5.838140201834775
$ws.Range("E7").Value = 10.68924794221167
$ws.Range("G7").Value = 50.13597244640479
$ws.Range("H7").Value = 19.30781233335457
$ws.Range("L7").Value = 9.988716433343523
$ws.Range("M7").Value = 17.0776135407258
$ws.Range("N7").Value = 20.27117540126042
$ws.Range("B8").Value = 20.30888808680814
$ws.Range("C8").Value = 8.920496647093884
$ws.Range("D8").Value = 5.97944339899334
$ws.Range("E8").Value = 10.68621184968114
$ws.Range("G8").Value = 50.65054919362019
$ws.Range("H8").Value = 19.30885222968424
$ws.Range("L8").Value = 9.972348469157199
$ws.Range("M8").Value = 17.16555051286444
$ws.Range("N8").Value = 20.1809079654957
$ws.Range("B9").Value = 21.34599499923041
$ws.Range("C9").Value = 10.01669271051243
$ws.Range("D9").Value = 6.261777323931987
$ws.Range("E9").Value = 10.68179989444929
$ws.Range("G9").Value = 51.77629745687496
$ws.Range("H9").Value = 19.3466705289309
$ws.Range("L9").Value = 9.951783516331911
$ws.Range("M9").Value = 17.36692692663159
$ws.Range("N9").Value = 20.01958281151399
$ws.Range("B10").Value = 22.10333673494963
$ws.Range("C10").Value = 10.744379215463
$ws.Range("D10").Value = 6.470083413680575
$ws.Range("E10").Value = 10.67950023837797
$ws.Range("G10").Value = 52.66836121984847
$ws.Range("H10").Value = 19.39630257265729
$ws.Range("L10").Value = 9.94369933882586
$ws.Range("M10").Value = 17.53167642847173
$ws.Range("N10").Value = 19.91055200522152
$ws.Range("B11").Value = 22.44514116673391
$ws.Range("C11").Value = 11.05828449572994
$ws.Range("D11").Value = 6.564566443590675
$ws.Range("E11").Value = 10.67865904801015
$ws.Range("G11").Value = 53.08685630714874
$ws.Range("H11").Value = 19.42361812605316
$ws.Range("L11").Value = 9.941544660621208
$ws.Range("M11").Value = 17.61007210940384
$ws.Range("N11").Value = 19.86299814216598
$ws.Range("B12").Value = 22.57404878949514
$ws.Range("C12").Value = 11.17467474103094
$ws.Range("D12").Value = 6.60026767093104
$ws.Range("E12").Value = 10.6783700198756
$ws.Range("G12").Value = 53.24702683601242
$ws.Range("H12").Value = 19.43464114093575
$ws.Range("L12").Value = 9.940947388682874
$ws.Range("M12").Value = 17.64023625549864
$ws.Range("N12").Value = 19.84528376232135
$ws.Range("B13").Value = 22.5463117162963
$ws.Range("C13").Value = 11.14971848191016
$ws.Range("D13").Value = 6.592582798517665
$ws.Range("E13").Value = 10.67843095397409
$ws.Range("G13").Value = 53.21245793341404
$ws.Range("H13").Value = 19.43223697030815
$ws.Range("L13").Value = 9.941066302941081
$ws.Range("M13").Value = 17.63371894156528
$ws.Range("N13").Value = 19.8490858417183
$ws.Range("B14").Value = 22.45575777357987
$ws.Range("C14").Value = 11.06790974977429
$ws.Range("D14").Value = 6.567505366769131
$ws.Range("E14").Value = 10.67863467776166
$ws.Range("G14").Value = 53.10000040521512
$ws.Range("H14").Value = 19.42451139827997
$ws.Range("L14").Value = 9.941491142583336
$ws.Range("M14").Value = 17.61254427431391
$ws.Range("N14").Value = 19.86153489710209
$ws.Range("B15").Value = 22.40021842559819
$ws.Range("C15").Value = 11.01747630564825
$ws.Range("D15").Value = 6.552133541508105
$ws.Range("E15").Value = 10.67876330917658
$ws.Range("G15").Value = 53.03133383094406
$ws.Range("H15").Value = 19.4198676411606
$ws.Range("L15").Value = 9.941779833985025
$ws.Range("M15").Value = 17.59963578059748
$ws.Range("N15").Value = 19.86919846755847
$ws.Range("B16").Value = 22.08093259386994
$ws.Range("C16").Value = 10.72351867779886
$ws.Range("D16").Value = 6.463899900802983
$ws.Range("E16").Value = 10.67955933957506
$ws.Range("G16").Value = 52.64125597206585
$ws.Range("H16").Value = 19.39461261181791
$ws.Range("L16").Value = 9.943870769473991
$ws.Range("M16").Value = 17.52662093215033
$ws.Range("N16").Value = 19.91370083929865
$ws.Range("B17").Value = 21.88426672510836
$ws.Range("C17").Value = 10.53878619363771
$ws.Range("D17").Value = 6.409673250521209
$ws.Range("E17").Value = 10.68010019690803
$ws.Range("G17").Value = 52.40511534837842
$ws.Range("H17").Value = 19.38033167183895
$ws.Range("L17").Value = 9.945543319450206
$ws.Range("M17").Value = 17.48269932781994
$ws.Range("N17").Value = 19.94152474482247
$ws.Range("B18").Value = 21.77090191609439
$ws.Range("C18").Value = 10.43092191030558
$ws.Range("D18").Value = 6.378459645243534
$ws.Range("E18").Value = 10.68043057042601
$ws.Range("G18").Value = 52.27049603093893
$ws.Range("H18").Value = 19.37256375784283
$ws.Range("L18").Value = 9.946648688138282
$ws.Range("M18").Value = 17.45776272073105
$ws.Range("N18").Value = 19.95772083189849
$ws.Range("B19").Value = 21.73248024929264
$ws.Range("C19").Value = 10.39412480596682
$ws.Range("D19").Value = 6.367888331861958
$ws.Range("E19").Value = 10.68054574046749
$ws.Range("G19").Value = 52.22512666157259
$ws.Range("H19").Value = 19.37001034523192
$ws.Range("L19").Value = 9.947047578399051
$ws.Range("M19").Value = 17.44937615140171
$ws.Range("N19").Value = 19.96323763776025
$ws.Range("B20").Value = 21.90522879903791
$ws.Range("C20").Value = 10.55861815458247
$ws.Range("D20").Value = 6.415448519083284
$ws.Range("E20").Value = 10.68004062531954
$ws.Range("G20").Value = 52.43012939060351
$ws.Range("H20").Value = 19.38180574471966
$ws.Range("L20").Value = 9.945350438362482
$ws.Range("M20").Value = 17.48734125735924
$ws.Range("N20").Value = 19.93854292417393
$ws.Range("B21").Value = 22.4823709904754
$ws.Range("C21").Value = 11.09200634213877
$ws.Range("D21").Value = 6.574873609814526
$ws.Range("E21").Value = 10.67857403783323
$ws.Range("G21").Value = 53.13298693905123
$ws.Range("H21").Value = 19.42676217271017
$ws.Range("L21").Value = 9.941360425511368
$ws.Range("M21").Value = 17.6187509842911
$ws.Range("N21").Value = 19.85787035569479
$ws.Range("B22").Value = 22.85644160246695
$ws.Range("C22").Value = 11.42615308736489
$ws.Range("D22").Value = 6.678599511037932
$ws.Range("E22").Value = 10.67778758096082
$ws.Range("G22").Value = 53.60216562040077
$ws.Range("H22").Value = 19.46010125572301
$ws.Range("L22").Value = 9.940027027223918
$ws.Range("M22").Value = 17.70740824768417
$ws.Range("N22").Value = 19.80685501831989
$ws.Range("B23").Value = 22.65712193277503
$ws.Range("C23").Value = 11.2491394489051
$ws.Range("D23").Value = 6.623293819017714
$ws.Range("E23").Value = 10.6781915709266
$ws.Range("G23").Value = 53.35090110036897
$ws.Range("H23").Value = 19.44194632674194
$ws.Range("L23").Value = 9.940622213543959
$ws.Range("M23").Value = 17.65984283322011
$ws.Range("N23").Value = 19.83392676277236
$ws.Range("B24").Value = 21.89575276627015
$ws.Range("C24").Value = 10.54965729150977
$ws.Range("D24").Value = 6.412837635591766
$ws.Range("E24").Value = 10.68006749712566
$ws.Range("G24").Value = 52.41881697815339
$ws.Range("H24").Value = 19.38113793828048
$ws.Range("L24").Value = 9.945437191981346
$ws.Range("M24").Value = 17.48524165997996
$ws.Range("N24").Value = 19.93989038304162
$ws.Range("B25").Value = 21.06563369175585
$ws.Range("C25").Value = 9.733672434755558
$ws.Range("D25").Value = 6.185080510748433
$ws.Range("E25").Value = 10.68282819862018
$ws.Range("G25").Value = 51.45987123525389
$ws.Range("H25").Value = 19.33260279846728
$ws.Range("L25").Value = 9.956112576004024
$ws.Range("M25").Value = 17.30943422722794
$ws.Range("N25").Value = 20.06155321106856
